# Append five new paragraphs (dropdown/theme-switch TODO items) right
# after the existing last paragraph of the document body ("... redo
# vertical nav bar root css variables I just lost (oof)"), matching the
# author's commit "Up to dropdown functionality to switch between theme,
# need to write backend".

$d = $word.ActiveDocument

# Locate the anchor paragraph by its distinctive text so the script is
# resilient to exact character-offset differences.
$anchorRange = $d.Content.Duplicate
$found = $anchorRange.Find.Execute(
    "redo vertical nav bar root css variables I just lost (oof)",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph text to insert after."
}

# Collapse to the end of the matched text, then grab the paragraph that
# contains it so we can insert exactly at that paragraph's end (i.e.
# right after its closing paragraph mark), regardless of where in the
# document it actually lives.
$anchorRange.Collapse(0)
$anchorParagraph = $anchorRange.Paragraphs(1)
$insertAt = $d.Range($anchorParagraph.Range.End, $anchorParagraph.Range.End)

$xmlPackage = @'
<?xml version="1.0" encoding="utf-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="284"/>
              </w:tabs>
              <w:spacing w:after="0"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="284"/>
              </w:tabs>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r>
              <w:tab/>
              <w:t>(todo2)</w:t>
            </w:r>
            <w:r>
              <w:tab/>
              <w:t>clean reset window/button</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="284"/>
              </w:tabs>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
              <w:t>Version number, by Jonius7</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="284"/>
              </w:tabs>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
              <w:t>run functions with reset=1 argument</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="284"/>
              </w:tabs>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r>
              <w:tab/>
              <w:t>(todo)</w:t>
            </w:r>
            <w:r>
              <w:tab/>
            </w:r>
            <w:r>
              <w:tab/>
              <w:t>Prompt for and download SteamFriendsPatcher automatically</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertAt.InsertXML($xmlPackage)
